$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "Contact" column (G) ---
$ws.Range("G1").Value = "Contact"

# --- Row 2: trim the task description, add contact name ---
$ws.Range("E2").Value = "Add user and login function"
$ws.Range("G2").Value = "Hà"

# --- Row 3: trim the task description, add contact name ---
$ws.Range("E3").Value = "Design ask and answer screen"
$ws.Range("G3").Value = "Bottle"

# --- Insert a new row after row 5 (old row6 "create a new category" shifts to 7,8,9) ---
$ws.Rows("6").Insert()

# --- Row 5 gains a second function row: "Show list categories" ---
$ws.Range("B5").Value = "On going"
$ws.Range("C5").Value = "Danh mục môn học"
$ws.Range("D5").Value = "Reddit"
$ws.Range("E5").Value = "Show list categories"
$ws.Range("F5").Value = "On going"

# Column C now holds longer text and got a manually fixed (non bestFit) width
$ws.Columns("C").ColumnWidth = 20.711495535714285

# --- New row 6: "Add category" task, Done ---
$ws.Range("E6").Value = "Add category"
$ws.Range("F6").Value = "Done"

# --- Move selection to A4, matching the saved workbook state ---
[void]$ws.Range("A4").Select()
